$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: merge the two runs "SAT Oct 28" + " 12:24:51 PDT 2017"
# into a single run's text.
# -----------------------------------------------------------------
$find = $d.Content.Find
$find.Text = "SAT Oct 28 12:24:51 PDT 2017"
$found = $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, "SAT Oct 28 12:24:51 PDT 2017", 2)
if (-not $found) {
    throw "Could not find 'SAT Oct 28 12:24:51 PDT 2017'"
}

# -----------------------------------------------------------------
# Change 2: after the "SAT Oct 28 ..." receipt block's
# "- CASH AND CLEARD" line, insert a whole new receipt entry
# (05/11/2017 MAMATHA CHICK IN -> "SAT Nov 04 10:57:32 PST 2017").
# -----------------------------------------------------------------

# Re-locate the (now merged) date-stamp paragraph, then search forward
# from there for the "- CASH AND CLEARD" line that belongs to that
# same receipt (there are several "- CASH AND CLEARD" lines earlier
# in the document, so we must anchor the search).
$dateFind = $d.Content.Find
$dateFind.Text = "SAT Oct 28 12:24:51 PDT 2017"
$dateFind.Execute() | Out-Null
$afterDate = $dateFind.Parent.End

$searchRange = $d.Range($afterDate, $d.Content.End)
$cashFind = $searchRange.Find
$cashFind.Text = "- CASH AND CLEARD"
$cashFound = $cashFind.Execute()
if (-not $cashFound) {
    throw "Could not find '- CASH AND CLEARD' after the SAT Oct 28 receipt"
}

# Collapse to the end of the found text (i.e. right before the
# paragraph mark that closes the "Amount Received mode ... - CASH AND
# CLEARD" paragraph), then insert the new paragraphs right there.
$searchRange.Collapse(0)
$insertPos = $searchRange.Start
$insertionRange = $d.Range($insertPos, $insertPos)

$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$font = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr>'
$fontB = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/></w:rPr>'

$pPr = "<w:pPr><w:pStyle w:val=`"PlainText`"/>$font</w:pPr>"
$pPrB = "<w:pPr><w:pStyle w:val=`"PlainText`"/>$fontB</w:pPr>"

$sb = New-Object System.Text.StringBuilder

# blank separator line
[void]$sb.Append("<w:p $W>$pPr</w:p>")

# date stamp of the new receipt
[void]$sb.Append("<w:p $W>$pPr<w:r>$font<w:t>SAT Nov 04</w:t></w:r><w:r>$font<w:t xml:space=`"preserve`"> 10:57:32 PST 2017</w:t></w:r></w:p>")

# Person Name
[void]$sb.Append("<w:p $W>$pPr<w:r>$font<w:t>Person Name</w:t></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/><w:t>- MAHADEVA</w:t></w:r></w:p>")

# Bill number
[void]$sb.Append("<w:p $W>$pPr<w:r>$font<w:t>Bill number</w:t></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/><w:t>- 1373</w:t></w:r></w:p>")

# dashed separator
[void]$sb.Append("<w:p $W>$pPr<w:r>$font<w:t>---------------------------------------------------------------</w:t></w:r></w:p>")

# Item Name
[void]$sb.Append("<w:p $W>$pPr<w:r>$font<w:t>Item Name</w:t></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/><w:t>- CARROT EVE</w:t></w:r></w:p>")

# Number of Pockets
[void]$sb.Append("<w:p $W>$pPr<w:r>$font<w:t>Number of Pockets</w:t></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/><w:t>- 1</w:t></w:r></w:p>")

# Number of KGs
[void]$sb.Append("<w:p $W>$pPr<w:r>$font<w:t>Number of KGs</w:t></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/><w:t>- 113</w:t></w:r></w:p>")

# Rate
[void]$sb.Append("<w:p $W>$pPr<w:r>$font<w:t>Rate</w:t></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/><w:t>- 42</w:t></w:r></w:p>")

# Transport & Miscellaneous
[void]$sb.Append("<w:p $W>$pPr<w:r>$font<w:t>Transport &amp; Miscellaneous</w:t></w:r><w:r>$font<w:tab/><w:t>- 10</w:t></w:r></w:p>")

# Total Price
[void]$sb.Append("<w:p $W>$pPr<w:r>$font<w:t>Total Price</w:t></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/></w:r><w:r>$font<w:tab/><w:t>- 4756.0</w:t></w:r></w:p>")

# Amount balance (bold)
[void]$sb.Append("<w:p $W>$pPrB<w:r>$fontB<w:t>Amount balance</w:t></w:r><w:r>$fontB<w:tab/></w:r><w:r>$fontB<w:tab/></w:r><w:r>$fontB<w:tab/><w:t>- 4756.0</w:t></w:r></w:p>")

# two trailing blank lines
[void]$sb.Append("<w:p $W>$pPr</w:p>")
[void]$sb.Append("<w:p $W>$pPr</w:p>")

$xml = $sb.ToString()
$insertionRange.InsertXML($xml) | Out-Null

Write-Host "Edit complete."
